# Re-generate the quadratic/linear experiment values (alpha_zero variant).
$wb = $excel.ActiveWorkbook

$ws_Restricciones_del_lider = $wb.Worksheets.Item(2)  # Restricciones_del_lider
$ws_Restricciones_del_lider.Range("A2").Value = '-16.45 + x_1 + x_2 + y_1 - 2y_2'

$ws_Restricciones_del_lider.Range("B2").NumberFormat = "@"
$ws_Restricciones_del_lider.Range("B2").Value = '-23.55'
$ws_Restricciones_del_lider.Range("B2").Style = "Normal"

$ws_Restricciones_del_lider.Range("D2").NumberFormat = "@"
$ws_Restricciones_del_lider.Range("D2").Value = '0.86'
$ws_Restricciones_del_lider.Range("D2").Style = "Normal"

$ws_Restricciones_del_follower = $wb.Worksheets.Item(3)  # Restricciones_del_follower
$ws_Restricciones_del_follower.Range("A2").Value = '18.85 - x_1 + 2y_2'

$ws_Restricciones_del_follower.Range("B2").NumberFormat = "@"
$ws_Restricciones_del_follower.Range("B2").Value = '-8.850000000000001'
$ws_Restricciones_del_follower.Range("B2").Style = "Normal"

$ws_Restricciones_del_follower.Range("D2").NumberFormat = "@"
$ws_Restricciones_del_follower.Range("D2").Value = '0.32'
$ws_Restricciones_del_follower.Range("D2").Style = "Normal"

$ws_Restricciones_del_follower.Range("E2").NumberFormat = "@"
$ws_Restricciones_del_follower.Range("E2").Value = '3.5999999999999996'
$ws_Restricciones_del_follower.Range("E2").Style = "Normal"

$ws_Restricciones_del_follower.Range("F2").NumberFormat = "@"
$ws_Restricciones_del_follower.Range("F2").Value = '7.7'
$ws_Restricciones_del_follower.Range("F2").Style = "Normal"

$ws_Restricciones_del_follower.Range("A3").Value = '20.0 - x_2 + 2y_2'

$ws_Restricciones_del_follower.Range("B3").NumberFormat = "@"
$ws_Restricciones_del_follower.Range("B3").Value = '-10.0'
$ws_Restricciones_del_follower.Range("B3").Style = "Normal"

$ws_Restricciones_del_follower.Range("D3").NumberFormat = "@"
$ws_Restricciones_del_follower.Range("D3").Value = '0.82'
$ws_Restricciones_del_follower.Range("D3").Style = "Normal"

$ws_Restricciones_del_follower.Range("E3").NumberFormat = "@"
$ws_Restricciones_del_follower.Range("E3").Value = '4.3'
$ws_Restricciones_del_follower.Range("E3").Style = "Normal"

$ws_Restricciones_del_follower.Range("F3").NumberFormat = "@"
$ws_Restricciones_del_follower.Range("F3").Value = '0.4'
$ws_Restricciones_del_follower.Range("F3").Style = "Normal"

$ws_Restricciones_del_follower.Range("A4").Value = '-5.1 - y_1'

$ws_Restricciones_del_follower.Range("B4").NumberFormat = "@"
$ws_Restricciones_del_follower.Range("B4").Value = '-4.9'
$ws_Restricciones_del_follower.Range("B4").Style = "Normal"

$ws_Restricciones_del_follower.Range("D4").NumberFormat = "@"
$ws_Restricciones_del_follower.Range("D4").Value = '0.26'
$ws_Restricciones_del_follower.Range("D4").Style = "Normal"

$ws_Restricciones_del_follower.Range("E4").NumberFormat = "@"
$ws_Restricciones_del_follower.Range("E4").Value = '7.800000000000001'
$ws_Restricciones_del_follower.Range("E4").Style = "Normal"

$ws_Restricciones_del_follower.Range("F4").NumberFormat = "@"
$ws_Restricciones_del_follower.Range("F4").Value = '4.2'
$ws_Restricciones_del_follower.Range("F4").Style = "Normal"

$ws_Restricciones_del_follower.Range("A5").Value = '5.100000000000001 + y_1'

$ws_Restricciones_del_follower.Range("B5").NumberFormat = "@"
$ws_Restricciones_del_follower.Range("B5").Value = '-25.1'
$ws_Restricciones_del_follower.Range("B5").Style = "Normal"

$ws_Restricciones_del_follower.Range("D5").NumberFormat = "@"
$ws_Restricciones_del_follower.Range("D5").Value = '0.85'
$ws_Restricciones_del_follower.Range("D5").Style = "Normal"

$ws_Restricciones_del_follower.Range("E5").NumberFormat = "@"
$ws_Restricciones_del_follower.Range("E5").Value = '1.9'
$ws_Restricciones_del_follower.Range("E5").Style = "Normal"

$ws_Restricciones_del_follower.Range("F5").NumberFormat = "@"
$ws_Restricciones_del_follower.Range("F5").Value = '0.3'
$ws_Restricciones_del_follower.Range("F5").Style = "Normal"

$ws_Restricciones_del_follower.Range("A6").Value = '-11.35 - y_2'

$ws_Restricciones_del_follower.Range("B6").NumberFormat = "@"
$ws_Restricciones_del_follower.Range("B6").Value = '-1.3499999999999996'
$ws_Restricciones_del_follower.Range("B6").Style = "Normal"

$ws_Restricciones_del_follower.Range("D6").NumberFormat = "@"
$ws_Restricciones_del_follower.Range("D6").Value = '0.7'
$ws_Restricciones_del_follower.Range("D6").Style = "Normal"

$ws_Restricciones_del_follower.Range("E6").NumberFormat = "@"
$ws_Restricciones_del_follower.Range("E6").Value = '6.1'
$ws_Restricciones_del_follower.Range("E6").Style = "Normal"

$ws_Restricciones_del_follower.Range("F6").NumberFormat = "@"
$ws_Restricciones_del_follower.Range("F6").Value = '5.8999999999999995'
$ws_Restricciones_del_follower.Range("F6").Style = "Normal"

$ws_Restricciones_del_follower.Range("A7").Value = '-48.65 + y_2'

$ws_Restricciones_del_follower.Range("B7").NumberFormat = "@"
$ws_Restricciones_del_follower.Range("B7").Value = '-28.65'
$ws_Restricciones_del_follower.Range("B7").Style = "Normal"

$ws_Restricciones_del_follower.Range("D7").NumberFormat = "@"
$ws_Restricciones_del_follower.Range("D7").Value = '0.51'
$ws_Restricciones_del_follower.Range("D7").Style = "Normal"

$ws_Restricciones_del_follower.Range("E7").NumberFormat = "@"
$ws_Restricciones_del_follower.Range("E7").Value = '1.5'
$ws_Restricciones_del_follower.Range("E7").Style = "Normal"

$ws_Restricciones_del_follower.Range("F7").NumberFormat = "@"
$ws_Restricciones_del_follower.Range("F7").Value = '8.5'
$ws_Restricciones_del_follower.Range("F7").Style = "Normal"

$ws_Punto_modificado = $wb.Worksheets.Item(4)  # Punto_modificado
$ws_Punto_modificado.Range("A2").NumberFormat = "@"
$ws_Punto_modificado.Range("A2").Value = '1.55'
$ws_Punto_modificado.Range("A2").Style = "Normal"

$ws_Punto_modificado.Range("B2").NumberFormat = "@"
$ws_Punto_modificado.Range("B2").Value = '2.7'
$ws_Punto_modificado.Range("B2").Style = "Normal"

$ws_Punto_modificado.Range("C2").NumberFormat = "@"
$ws_Punto_modificado.Range("C2").Value = '-5.1'
$ws_Punto_modificado.Range("C2").Style = "Normal"

$ws_Punto_modificado.Range("D2").NumberFormat = "@"
$ws_Punto_modificado.Range("D2").Value = '-8.65'
$ws_Punto_modificado.Range("D2").Style = "Normal"

$ws_Vector_bf = $wb.Worksheets.Item(5)  # Vector_bf
$ws_Vector_bf.Range("A2").NumberFormat = "@"
$ws_Vector_bf.Range("A2").Value = '-27.29'
$ws_Vector_bf.Range("A2").Style = "Normal"

$ws_Vector_bf.Range("A3").NumberFormat = "@"
$ws_Vector_bf.Range("A3").Value = '-19.39'
$ws_Vector_bf.Range("A3").Style = "Normal"

$ws_Vector_BF = $wb.Worksheets.Item(6)  # Vector_BF
$ws_Vector_BF.Range("A2").NumberFormat = "@"
$ws_Vector_BF.Range("A2").Value = '0.7399999999999998'
$ws_Vector_BF.Range("A2").Style = "Normal"

$ws_Vector_BF.Range("A3").NumberFormat = "@"
$ws_Vector_BF.Range("A3").Value = '1.44'
$ws_Vector_BF.Range("A3").Style = "Normal"

$ws_Vector_BF.Range("A4").NumberFormat = "@"
$ws_Vector_BF.Range("A4").Value = '8.040000000000001'
$ws_Vector_BF.Range("A4").Style = "Normal"

$ws_Vector_BF.Range("A5").NumberFormat = "@"
$ws_Vector_BF.Range("A5").Value = '-6.4799999999999995'
$ws_Vector_BF.Range("A5").Style = "Normal"

